$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 updates
$ws.Range("G3").Value = 1.95
$ws.Range("H3").Value = 3.3
$ws.Range("I3").Value = 4.2
$ws.Range("L3").Value = 5
$ws.Range("AB3").Value = 41
$ws.Range("AD3").Value = 6.5
$ws.Range("AT3").Value = 9.5
$ws.Range("AY3").Value = 101
$ws.Range("AZ3").Value = 151
$ws.Range("BA3").Value = 401

# Row 4 updates
$ws.Range("G4").Value = 2.6
$ws.Range("I4").Value = 3.1
$ws.Range("M4").Value = 1.14
$ws.Range("N4").Value = 5.5
$ws.Range("S4").Value = 1.75
$ws.Range("T4").Value = 2.05
$ws.Range("Z4").Value = 26
$ws.Range("AE4").Value = 23
$ws.Range("AN4").Value = 17
$ws.Range("AY4").Value = 81

# Row 5 updates
$ws.Range("M5").Value = 1.04
$ws.Range("N5").Value = 13
